$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range("G2").Value = 27.24688066666667
$ws.Range("H2").Value = 81.74064200000001
$ws.Range("I2").Value = 0.9344428125547477
$ws.Range("J2").Value = 0.9344428125547476
$ws.Range("M2").Value = 12.39940933333333
$ws.Range("N2").Value = 37.198228
$ws.Range("O2").Value = 0.6889801160127385
$ws.Range("P2").Value = 0.6889801160127385
$ws.Range("Q2").Value = 337.8452264424863
$ws.Range("R2").Value = 3040.607037982376
$ws.Range("S2").Value = 0.6438125174012397
$ws.Range("T2").Value = 0.6438125174012396
$ws.Range("G3").Value = 27.24688066666667
$ws.Range("H3").Value = 81.74064200000001
$ws.Range("I3").Value = 0.9344428125547477
$ws.Range("J3").Value = 0.9344428125547476
$ws.Range("O3").Value = 0.1923952864254561
$ws.Range("P3").Value = 0.1923952864254561
$ws.Range("Q3").Value = 94.34209725099446
$ws.Range("R3").Value = 849.0788752589501
$ws.Range("S3").Value = 0.1797823925696795
$ws.Range("T3").Value = 0.1797823925696794
$ws.Range("G4").Value = 27.24688066666667
$ws.Range("H4").Value = 81.74064200000001
$ws.Range("I4").Value = 0.9344428125547477
$ws.Range("J4").Value = 0.9344428125547476
$ws.Range("M4").Value = 2.134858333333333
$ws.Range("N4").Value = 6.404574999999999
$ws.Range("O4").Value = 0.1186245975618055
$ws.Range("P4").Value = 0.1186245975618055
$ws.Range("Q4").Value = 58.16823024857223
$ws.Range("R4").Value = 523.51407223715
$ws.Range("S4").Value = 0.1108479025838286
$ws.Range("T4").Value = 0.1108479025838286
$ws.Range("H5").Value = 0.862752
$ws.Range("I5").Value = 0.009862809805399298
$ws.Range("J5").Value = 0.009862809805399296
$ws.Range("M5").Value = 12.39940933333333
$ws.Range("N5").Value = 37.198228
$ws.Range("O5").Value = 0.6889801160127385
$ws.Range("P5").Value = 0.6889801160127385
$ws.Range("Q5").Value = 3.565871733717333
$ws.Range("R5").Value = 32.092845603456
$ws.Range("S5").Value = 0.006795279843935582
$ws.Range("T5").Value = 0.006795279843935582
$ws.Range("H6").Value = 0.862752
$ws.Range("I6").Value = 0.009862809805399298
$ws.Range("J6").Value = 0.009862809805399296
$ws.Range("O6").Value = 0.1923952864254561
$ws.Range("P6").Value = 0.1923952864254561
$ws.Range("R6").Value = 8.9618148312
$ws.Range("S6").Value = 0.001897558117469595
$ws.Range("T6").Value = 0.001897558117469594
$ws.Range("H7").Value = 0.862752
$ws.Range("I7").Value = 0.009862809805399298
$ws.Range("J7").Value = 0.009862809805399296
$ws.Range("M7").Value = 2.134858333333333
$ws.Range("N7").Value = 6.404574999999999
$ws.Range("O7").Value = 0.1186245975618055
$ws.Range("P7").Value = 0.1186245975618055
$ws.Range("Q7").Value = 0.6139510989333333
$ws.Range("R7").Value = 5.525559890399999
$ws.Range("S7").Value = 0.001169971843994121
$ws.Range("T7").Value = 0.00116997184399412
$ws.Range("G8").Value = 0.2337766666666667
$ws.Range("H8").Value = 0.70133
$ws.Range("I8").Value = 0.008017465506681745
$ws.Range("J8").Value = 0.008017465506681745
$ws.Range("M8").Value = 12.39940933333333
$ws.Range("N8").Value = 37.198228
$ws.Range("O8").Value = 0.6889801160127385
$ws.Range("P8").Value = 0.6889801160127385
$ws.Range("Q8").Value = 2.898692582582222
$ws.Range("R8").Value = 26.08823324324
$ws.Range("S8").Value = 0.005523874314921718
$ws.Range("T8").Value = 0.005523874314921718
$ws.Range("G9").Value = 0.2337766666666667
$ws.Range("H9").Value = 0.70133
$ws.Range("I9").Value = 0.008017465506681745
$ws.Range("J9").Value = 0.008017465506681745
$ws.Range("O9").Value = 0.1923952864254561
$ws.Range("P9").Value = 0.1923952864254561
$ws.Range("Q9").Value = 0.8094497601944445
$ws.Range("R9").Value = 7.28504784175
$ws.Range("S9").Value = 0.001542522572564249
$ws.Range("T9").Value = 0.001542522572564249
$ws.Range("G10").Value = 0.2337766666666667
$ws.Range("H10").Value = 0.70133
$ws.Range("I10").Value = 0.008017465506681745
$ws.Range("J10").Value = 0.008017465506681745
$ws.Range("M10").Value = 2.134858333333333
$ws.Range("N10").Value = 6.404574999999999
$ws.Range("O10").Value = 0.1186245975618055
$ws.Range("P10").Value = 0.1186245975618055
$ws.Range("Q10").Value = 0.4990800649722222
$ws.Range("R10").Value = 4.491720584749999
$ws.Range("S10").Value = 0.0009510686191957789
$ws.Range("T10").Value = 0.0009510686191957786
$ws.Range("G11").Value = 1.390183666666666
$ws.Range("H11").Value = 4.170551
$ws.Range("I11").Value = 0.04767691213317134
$ws.Range("J11").Value = 0.04767691213317134
$ws.Range("M11").Value = 12.39940933333333
$ws.Range("N11").Value = 37.198228
$ws.Range("O11").Value = 0.6889801160127385
$ws.Range("P11").Value = 0.6889801160127385
$ws.Range("Q11").Value = 17.23745633151422
$ws.Range("R11").Value = 155.137106983628
$ws.Range("S11").Value = 0.03284844445264153
$ws.Range("T11").Value = 0.03284844445264153
$ws.Range("G12").Value = 1.390183666666666
$ws.Range("H12").Value = 4.170551
$ws.Range("I12").Value = 0.04767691213317134
$ws.Range("J12").Value = 0.04767691213317134
$ws.Range("O12").Value = 0.1923952864254561
$ws.Range("P12").Value = 0.1923952864254561
$ws.Range("Q12").Value = 4.813499360969444
$ws.Range("R12").Value = 43.321494248725
$ws.Range("S12").Value = 0.009172813165742802
$ws.Range("T12").Value = 0.0091728131657428
$ws.Range("G13").Value = 1.390183666666666
$ws.Range("H13").Value = 4.170551
$ws.Range("I13").Value = 0.04767691213317134
$ws.Range("J13").Value = 0.04767691213317134
$ws.Range("M13").Value = 2.134858333333333
$ws.Range("N13").Value = 6.404574999999999
$ws.Range("O13").Value = 0.1186245975618055
$ws.Range("P13").Value = 0.1186245975618055
$ws.Range("Q13").Value = 2.967845185647222
$ws.Range("R13").Value = 26.71060667082499
$ws.Range("S13").Value = 0.00565565451478701
$ws.Range("T13").Value = 0.005655654514787009
